# Fixed title sizing bug, Added KRI 10, 222, 273
#
# The "Title 3" placeholder shape on slide layouts 2, 3 and 4 (the BLANK /
# 1_Blank KRI section-header layouts) had its horizontal position/width
# adjusted so the title text box better fits the available space.
#
# PowerPoint's Shape.Left/.Width are expressed in points; the point values
# below were chosen so that converting them back to EMU (PowerPoint's
# native unit, 12700 EMU per point) reproduces the exact target EMU
# offsets/extents from the authoritative OOXML.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# layout index -> new Left (pt), new Width (pt)  [Top/Height are unchanged]
$changes = @{
    2 = @{ left = 42.7826;   width = 873.39133 }   # off x=543339  ext cx=11092070
    3 = @{ left = 42.7826;   width = 874.4348  }   # off x=543339  ext cx=11105322
    4 = @{ left = 43.8261;   width = 873.39133 }   # off x=556591  ext cx=11092070
}

foreach ($layoutIndex in $changes.Keys) {
    $cl = $m.CustomLayouts.Item($layoutIndex)

    foreach ($shp in $cl.Shapes) {
        if ($shp.Name -eq "Title 3") {
            $shp.Left  = $changes[$layoutIndex].left
            $shp.Width = $changes[$layoutIndex].width
        }
    }
}
